$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Headers: default (Item 1) -> header2.xml, first-page (Item 2) -> header1.xml.
# Both contain the "BTec_Logo-Orange" inline picture whose name swaps
# from image1.jpg to image2.jpg.
for ($hIdx = 1; $hIdx -le 2; $hIdx++) {
    $hdr = $sec.Headers.Item($hIdx)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

# Footers: default (Item 1) -> footer2.xml, first-page (Item 2) -> footer1.xml.
# Both contain the Pearson Edexcel logo inline picture whose name swaps
# from image2.png to image1.png.
for ($fIdx = 1; $fIdx -le 2; $fIdx++) {
    $ftr = $sec.Footers.Item($fIdx)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
